$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXHome")

# Insert a new column before column H ("Price") for the new "Side" field.
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column H.
$ws.Cells.Item(1, 8).Value = "Side"

# Fill in row 13 (previously blank test data row for QA_TestCase_Auto_NitroX_015).
$ws.Cells.Item(13, 2).Value = "Spot"
$ws.Cells.Item(13, 3).Value = "Trader01@Tinyex"
$ws.Cells.Item(13, 4).Value = "ETH"
$ws.Cells.Item(13, 5).Value = "USDT"
$ws.Cells.Item(13, 11).Value = 1

# Add new row 14 for QA_TestCase_Auto_NitroX_015_01, copying row 13's formatting.
$ws.Range("A13:P13").Copy($ws.Range("A14:P14"))
$ws.Cells.Item(14, 1).Value = "QA_TestCase_Auto_NitroX_015_01"
$ws.Rows.Item(14).RowHeight = 29.5

# Widen column A to fit the new longer test case name.
$ws.Columns.Item(1).ColumnWidth = 32.5

# Update selection / view state to match the saved workbook.
$ws.Range("E13").Select()

$wb.Windows.Item(1).Height = 10420
